# "se cambio estado en la CP-003"
# Update the CP-003 row (row 6) of the test log:
#  - Fecha (date) moves from 2024-10-31 to 2024-11-01
#  - Estado goes from "Fallido" to "Resuelto"
#  - Comentarios goes from "Se creo turno" to "Se mostro error turno invalido"
# Column F is widened to fit the new, longer comment text, and the
# active selection is left on F7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = 45597
$ws.Range("E6").Value = "Resuelto"
$ws.Range("F6").Value = "Se mostro error turno invalido"

# Re-fit column F (Comentarios) now that it holds the longer text.
$ws.Columns(6).ColumnWidth = 23.6

$ws.Range("F7").Select() | Out-Null
